$wb = $excel.ActiveWorkbook

# --- Metadata sheet: update URL, Version, Date, Publisher -----------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/plan-employee-code"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Elements sheet: update Fixed Value URL for Extension.url, and clear --
# --- the mis-placed ele-1/ext-1 constraint text on the Extension row ------
$elements = $wb.Worksheets.Item("Elements")
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/plan-employee-code"
$elements.Range("AI2").Value = ""
